$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: mark Background Processing as TRUE
$ws.Range("G2").Value = $true

# Remove row 3 entirely (it was a duplicate test-plan row)
$ws.Rows("3:3").Delete()

# Reflect the selection state captured in the saved workbook (entire row 2 selected)
$ws.Range("A2:XFD2").Select()
